# Apply the changes described by the diff:
# 1. Change D9 value from "annee_id" to "semestre_id"
# 2. Change C10 value from "annee_id" to "semestre_id"
# 3. Change the active cell / selection on Sheet1 from D11 to C11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values (shared string change: annee_id -> semestre_id)
$ws.Range("D9").Value = "semestre_id"
$ws.Range("C10").Value = "semestre_id"

# Update the selection / active cell
$ws.Range("C11").Select()
